$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13, columns A-D; row13's E value is removed entirely.
# Values given as strings then cast to [double] to safely handle scientific
# notation literals (e.g. "2.137739203028708e-06").
$data = @(
    @("0.1387899531097405", "0.03914180814872054", "3.545823754038231", "0.0003913878810184778", "1.466"),
    @("0.1321767733134403", "0.03918708777761198", "3.372967495404299", "0.0007436273065394205", "1.454"),
    @("0.139010761023405",  "0.02932749147019469", "4.739947198166628", "2.137739203028708e-06", "1.009"),
    @("0.1789237596779445", "0.03848168086892773", "4.649582753086483", "3.326072131377621e-06", "1.482"),
    @("0.5896458485383199", "0.03375276505388262", "17.46955686732671", "2.44385415150239e-68",  "1.068"),
    @("0.1690423234800643", "0.03391660016577125", "4.984058621850382", "6.226420173367331e-07", "1.119"),
    @("0.1598505233033204", "0.0334548220971973",  "4.778101131098585", "1.769583782460926e-06", "1.018"),
    @("0.1308310092372296", "0.03312073277142814", "3.950124236082482", "7.811063847392322e-05", "1.109"),
    @("0.09441272582165894","0.02742322653543831", "3.442801513514534", "0.0005757216924978548", "1.015"),
    @("0.1348725950536856", "0.03113009742329271", "4.332546513419159", "1.473945101251338e-05", "1.109"),
    @("0.1454387032560523", "0.03085162435369979", "4.71413438685315",  "2.427404268360595e-06", "1.053"),
    @("-1.427338443907308", "0.03392686863826552", "-42.07103399744466", "0", $null)
)

$rowIndex = 2
foreach ($rowVals in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = [double]$rowVals[0]
    $ws.Cells.Item($rowIndex, 2).Value = [double]$rowVals[1]
    $ws.Cells.Item($rowIndex, 3).Value = [double]$rowVals[2]
    $ws.Cells.Item($rowIndex, 4).Value = [double]$rowVals[3]
    if ($null -eq $rowVals[4]) {
        $ws.Cells.Item($rowIndex, 5).ClearContents()
    } else {
        $ws.Cells.Item($rowIndex, 5).Value = [double]$rowVals[4]
    }
    $rowIndex++
}

# Remove old rows 14 and 15 entirely -> used range shrinks to A1:E13
$ws.Range("A14:E15").Delete()
